$wb = $excel.ActiveWorkbook

# --- Config sheet: update the DB engine value from "MySQL" to "MySQLi" ---
$wsConfig = $wb.Worksheets.Item("Config")
$wsConfig.Range("A2").Value = "MySQLi"

# --- ZeroCode sheet: same engine value update, plus re-point the remembered selection ---
$wsZeroCode = $wb.Worksheets.Item("ZeroCode")
$wsZeroCode.Range("E2").Value = "MySQLi"
[void]$wsZeroCode.Range("E2").Select()

# --- Make Config the active/selected sheet & cell (was Paises before) ---
[void]$wsConfig.Activate()
[void]$wsConfig.Range("A2").Select()
